$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfa"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4732933333333333
$ws.Range("H2").Value = 1.41988
$ws.Range("I2").Value = 0.6351867741147567
$ws.Range("J2").Value = 0.6351867741147565
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.232494
$ws.Range("N2").Value = 6.697482000000001
$ws.Range("O2").Value = 0.2983035867032896
$ws.Range("P2").Value = 0.2983035867032895
$ws.Range("Q2").Value = 1.056624526906667
$ws.Range("R2").Value = 9.509620742160001
$ws.Range("S2").Value = 0.1894784929449241
$ws.Range("T2").Value = 0.1894784929449241

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfa"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4732933333333333
$ws.Range("H3").Value = 1.41988
$ws.Range("I3").Value = 0.6351867741147567
$ws.Range("J3").Value = 0.6351867741147565
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.575715666666667
$ws.Range("N3").Value = 10.727147
$ws.Range("O3").Value = 0.4777835050834676
$ws.Range("P3").Value = 0.4777835050834675
$ws.Range("Q3").Value = 1.692362386928889
$ws.Range("R3").Value = 15.23126148236
$ws.Range("S3").Value = 0.3034817633192092
$ws.Range("T3").Value = 0.3034817633192091

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfa"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4732933333333333
$ws.Range("H4").Value = 1.41988
$ws.Range("I4").Value = 0.6351867741147567
$ws.Range("J4").Value = 0.6351867741147565
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.675756666666667
$ws.Range("N4").Value = 5.02727
$ws.Range("O4").Value = 0.2239129082132429
$ws.Range("P4").Value = 0.2239129082132428
$ws.Range("Q4").Value = 0.7931244586222221
$ws.Range("R4").Value = 7.1381201276
$ws.Range("S4").Value = 0.1422265178506233
$ws.Range("T4").Value = 0.1422265178506233

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Tgfa"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2718313333333334
$ws.Range("H5").Value = 0.815494
$ws.Range("I5").Value = 0.3648132258852434
$ws.Range("J5").Value = 0.3648132258852433
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.232494
$ws.Range("N5").Value = 6.697482000000001
$ws.Range("O5").Value = 0.2983035867032896
$ws.Range("P5").Value = 0.2983035867032895
$ws.Range("Q5").Value = 0.6068618206786669
$ws.Range("R5").Value = 5.461756386108001
$ws.Range("S5").Value = 0.1088250937583655
$ws.Range("T5").Value = 0.1088250937583654

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Tgfa"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2718313333333334
$ws.Range("H6").Value = 0.815494
$ws.Range("I6").Value = 0.3648132258852434
$ws.Range("J6").Value = 0.3648132258852433
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.575715666666667
$ws.Range("N6").Value = 10.727147
$ws.Range("O6").Value = 0.4777835050834676
$ws.Range("P6").Value = 0.4777835050834675
$ws.Range("Q6").Value = 0.9719915572908892
$ws.Range("R6").Value = 8.747924015618
$ws.Range("S6").Value = 0.1743017417642584
$ws.Range("T6").Value = 0.1743017417642583

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Tgfa"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2718313333333334
$ws.Range("H7").Value = 0.815494
$ws.Range("I7").Value = 0.3648132258852434
$ws.Range("J7").Value = 0.3648132258852433
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.675756666666667
$ws.Range("N7").Value = 5.02727
$ws.Range("O7").Value = 0.2239129082132429
$ws.Range("P7").Value = 0.2239129082132428
$ws.Range("Q7").Value = 0.4555231690422222
$ws.Range("R7").Value = 4.09970852138
$ws.Range("S7").Value = 0.08168639036261954
$ws.Range("T7").Value = 0.08168639036261952
